$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new columns O1:V1 (Tipo..Fecha Ultimo Mantenimiento) ---
# Copy N1 formatting (bold, centered, bordered) onto the new header cells first
$ws.Range("N1").Copy()
$ws.Range("O1:V1").PasteSpecial(-4122)
$ws.Range("O1").Value = "Tipo"
$ws.Range("P1").Value = "Modelo"
$ws.Range("Q1").Value = "Número Serie"
$ws.Range("R1").Value = "Gabinete"
$ws.Range("S1").Value = "Total Puertos"
$ws.Range("T1").Value = "PoE"
$ws.Range("U1").Value = "Fecha Instalación"
$ws.Range("V1").Value = "Fecha Último Mantenimiento"

# --- New switch rows 5-7 (SW-004, SW-005, SW-006) ---
# Values are written column-by-column to match original authoring order
$ws.Range("A5").Value = "SW-004"
$ws.Range("A6").Value = "SW-005"
$ws.Range("A7").Value = "SW-006"
$ws.Range("B5").Value = "UBI-006"
$ws.Range("B6").Value = "UBI-007"
$ws.Range("B7").Value = "UBI-004"
$ws.Range("D5").Value = "TP-Link"
$ws.Range("D6").Value = "Cisco"
$ws.Range("D7").Value = "Ubiquiti"
$ws.Range("H5").Value = 5
$ws.Range("H6").Value = 15
$ws.Range("H7").Value = 8
$ws.Range("I5").Value = 3
$ws.Range("I6").Value = 9
$ws.Range("I7").Value = 0
$ws.Range("K5").Value = "Funcionando"
$ws.Range("K6").Value = "Funcionando"
$ws.Range("K7").Value = "Funcionando"
$ws.Range("N5").Value = "Switch campus Pucón"
$ws.Range("N6").Value = "Switch principal CFT Prat"
$ws.Range("N7").Value = "Sin puertos disponibles - considerar ampliación"
$ws.Range("O5").Value = "Switch PoE 8 puertos"
$ws.Range("O6").Value = "Switch PoE 24 puertos"
$ws.Range("O7").Value = "Switch PoE 8 puertos"
$ws.Range("P5").Value = "TL-SG1008P"
$ws.Range("P6").Value = "SG350-28P"
$ws.Range("P7").Value = "US-8-150W"
$ws.Range("Q5").Value = "SN-SW004"
$ws.Range("Q6").Value = "SN-SW005"
$ws.Range("Q7").Value = "SN-SW006"
$ws.Range("R5").Value = "GAB-004"
$ws.Range("R6").Value = "GAB-005"
$ws.Range("R7").Value = "GAB-006"
$ws.Range("S5").Value = 8
$ws.Range("S6").Value = 24
$ws.Range("S7").Value = 8
$ws.Range("T5").Value = "Sí"
$ws.Range("T6").Value = "Sí"
$ws.Range("T7").Value = "Sí"

# Installation/maintenance dates for new rows must stay as plain text (not auto-parsed as dates).
# Pre-format as Text, enter the values, then clear the format back off so the cells keep the
# workbook default style (matches how the source data was authored). Touch only the cells that
# actually get a value - row 5 has no "Fecha Ultimo Mantenimiento" (V5 must stay blank/absent).
$ws.Range("U5").NumberFormat = "@"
$ws.Range("U6").NumberFormat = "@"
$ws.Range("U7").NumberFormat = "@"
$ws.Range("V6").NumberFormat = "@"
$ws.Range("V7").NumberFormat = "@"
$ws.Range("U5").Value = "2024-08-05"
$ws.Range("U6").Value = "2024-05-12"
$ws.Range("U7").Value = "2024-06-20"
$ws.Range("V6").Value = "2025-09-18"
$ws.Range("V7").Value = "2025-10-10"
$ws.Range("U5").ClearFormats()
$ws.Range("U6").ClearFormats()
$ws.Range("U7").ClearFormats()
$ws.Range("V6").ClearFormats()
$ws.Range("V7").ClearFormats()
